$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 17-20 (old rows), since the new dimension is A1:C16
$ws.Range("A17:C20").EntireRow.Delete() | Out-Null

# Clear C2 (becomes an empty inline string cell)
$ws.Range("C2").Value = ""

# Update C3:C16 with new transcript text
$ws.Range("C3").Value = 'ライフとライフの映画のどっちがいいか っていうのいいか聞いていきたいと思います 僕はナイフを持って行く俺は 俺は ライターを持って行きます 僕も ライターの方がいいと思います 僕も会いたいと思います'
$ws.Range("C4").Value = 'じゃあじゃあじゃあ ついてきてきたいと思います 僕の 辛い理由としてはナイフは最悪なくてもマイなんとかなると思うんですけど ライターの場合は夜とかに夜が寒くて寒いので逆にあったかくできるものがないと逆に夜を過ごすことができないと思います 相手を選んだ理由が まず ナイフと比較した時に入ります ってなんか正直 そこら辺の石とかをなんかいい感じなかったら 砕いてない服の代役にはなるかな'
$ws.Range("C5").Value = '音村 天体 薬が効かない原始的に火の起こし方をしない限り できないことの方が便利かなナイフ 刃こぼれしたらもう使えないと焦げるかわかんないけどライターはオイルがなくなっても動物とか最悪自分の油とかで燃料にしてるんじゃないかというよりかはまだ 長く使い続けるかしました'
$ws.Range("C6").Value = 'ナイフって銀行に行ったらもしかしたら その真ん中に付いてたりしても使えるかもしれん すけど ライターって長い 過ぎたらもうほ 水着で使えないとかそう考えると思い灯 松任 ライター 火を起こすと やっぱり通すとかも買ってこないのます 松明とかも作れるし 本当に最悪な状況になった場合にあの山全体を燃やしてその煙で休場 なんだり ニュースになることもあるかもしれないよ 想定してま ライターでいちかばちか では'
$ws.Range("C7").Value = 'に関しては気を10時に そんなみんな と思ってるとかじゃなくて あの木の棒 2本と伝え 一本で そのこ しこしこするやり方ですぐにしよう 起こせるんで俺は 俺はできるけどできる人とできない人は分かるんだ あればできない人もいるんですよ 多分 燃料が切れたら多分出たらそれを解決するんですけどそのそんな あの その さっきの子を起こす アイデアはない 人が燃料を補充できるとは思えないんで 燃料入れたらその原始的な方法を使ってことはその心理学的に多分しよう 小杉も食べないと思うんでナイフを持ってたら確実に 最初からその子を起こすとどう思うんでなんで 僕はない'
$ws.Range("C8").Value = 'なんで僕はない'
$ws.Range("C9").Value = '無人島で動物とかいるのある人いますかじゃなくて動物にそのナイフで応戦とか狩りができるって言うので 笑'
$ws.Range("C10").Value = 'ライター 竹田 きついですよね そんな綺麗 女ライターに立ってる人なんて そんな 木の棒でまずあの団 撮ることしか頭にないんでそんななんかね 動物かどうなって 私を起こしやすい時間的とか心理的な負荷 楽 一日一回 しょう子のボートがなんとかあとないとか 準備するの めちゃくちゃ めんどくさい しちゃうと思う その方が効率いいんじゃないかなと すぐに使用するんで 緊急時とかに虫とか動物とかが襲われた時に使ってるとか ドア 照明'
$ws.Range("C11").Value = '照明があり 暗い 足元で彼の意見とディスカバリーチャンネル 見たい番組表を持って 正直な感想 述べると子を起こす 絵を見せたほうが 視聴者が喜ぶからあたかも そういう演出をしているだけであって 正直ナイフ 使ってる時点で どうなんかなって そこはちょっと微妙なんじゃないか 大富豪'
$ws.Range("C12").Value = 'ナイフ持ってたら それなんかちょっとそれ ライター 使わなくてもない 2日って体臭 ちゃう みたいな 思っちゃうんですけど今は食料調達として まず魚とかどう しなくても 見劣り玉 出てきのみを撮る時にナイフがないと取れないと無理だよ そこはそんなに頼ってるやったらわかるやんで 僕みたいに肉体があればいけると思いがあれば いらないと思うけどその通りだと何分ですか'
$ws.Range("C13").Value = '神ちゃんねるとかで結局 ナンダカンダ ってますかわからない 足りないとかなんか 道具とか釣竿とか作るになることしないと思うんですけどそこは どうでもいいし でも削れるくらいですか意思と意志を当てて怒らせて張ってたりとか それこそ恐竜したガスとか マフィンとかをバット 歌詞'
$ws.Range("C14").Value = '走ったらすぐに尖ったものとか入りなもの作れると思うよ そう考えた時にそういった道具を使ったら一応 ナイフの代役としても使えるんで そっちの方が良かったらまだ ライター持ってった方が楽なんじゃないかな 人はいますか'
$ws.Range("C15").Value = '使えないんでそんなライター 落ちた子に行けなければ まとめ の方に入れるどうでしょうかね 俺は 対処法に言った通り やったから さもうライターの方は ミントに落ちてないこと 誰かが私のことを起こすことが簡単であるため えんため またやっちゃんは'
$ws.Range("C16").Value = 'やっちゃんは終わった雛様簡単元 無人島に置いて こっちでも使えないってことあり また楽に日に火を起こしてる 止められた葉であり ナイフは木の実をとることや 万が一に備えた時の時とかに使える使うため ガンマナイフ波みたいな感じで まぁ 人それぞれだと思いますが 自宅だと思います'

Write-Host "Done. UsedRange: " $ws.UsedRange.Address()
